# Ajout d'images sur le powerpoint
#
# - Removes the red "Insérer ..." placeholder textboxes and the leftover
#   "Tableau 4" sample table / "Image 1" sample pictures used as visual
#   guides on slides 1, 4, 5 and 6.
# - Adds a blank themed "photo" placeholder shape (style-referenced
#   rectangle, no fill) on slides 4, 5, 6 and 7 where users are meant to
#   drop their own picture, mirroring the one that already exists on
#   slide 1 ("photoBatiment").

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1: drop the red instructional textbox "ZoneTexte 6"
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item("ZoneTexte 6").Delete()

# Template shape used to replicate the themed "photo" placeholder
# (lnRef/fillRef/effectRef/fontRef style, noFill, centered empty text)
# onto the other slides via Copy/Paste.
$photoTemplate = $s1.Shapes.Item("photoBatiment")

# ---------------------------------------------------------------------
# Slide 4: drop the sample "Tableau 4" table and the red "Insérer les
# photos" textbox, then add the new "photo" placeholder shape.
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item("Tableau 4").Delete()
$s4.Shapes.Item("ZoneTexte 3").Delete()

$photoTemplate.Copy()
$newShape = $s4.Shapes.Paste().Item(1)
$newShape.Name = "photo"
$newShape.Left = 506.8363952636719
$newShape.Top = 79.27000427246094
$newShape.Width = 433.50023622047246
$newShape.Height = 332.63702392578125

# ---------------------------------------------------------------------
# Slide 5: drop the sample "Image 1" picture and the red "Insérer les
# photos" textbox, then add the new "photo" placeholder shape.
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item("Image 1").Delete()
$s5.Shapes.Item("ZoneTexte 4").Delete()

$photoTemplate.Copy()
$newShape = $s5.Shapes.Paste().Item(1)
$newShape.Name = "photo"
$newShape.Left = 506.8363952636719
$newShape.Top = 79.27000427246094
$newShape.Width = 433.50023622047246
$newShape.Height = 332.63702392578125

# ---------------------------------------------------------------------
# Slide 6: drop the sample "Image 1" picture and the red "Insérer les
# photos" textbox, then add the new "photo" placeholder shape.
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item("Image 1").Delete()
$s6.Shapes.Item("ZoneTexte 4").Delete()

$photoTemplate.Copy()
$newShape = $s6.Shapes.Paste().Item(1)
$newShape.Name = "photo"
$newShape.Left = 506.8363952636719
$newShape.Top = 79.27000427246094
$newShape.Width = 433.50023622047246
$newShape.Height = 332.63702392578125

# ---------------------------------------------------------------------
# Slide 7: just add the new "photo" placeholder shape (different size,
# aligned next to "tableauPreconisations").
# ---------------------------------------------------------------------
$s7 = $p.Slides.Item(7)

$photoTemplate.Copy()
$newShape = $s7.Shapes.Paste().Item(1)
$newShape.Name = "photo"
$newShape.Left = 432.83734130859375
$newShape.Top = 92.83464566929133
$newShape.Width = 497.45294189453125
$newShape.Height = 433.1021423339844
